$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A, B and E in rows 103-111 hold numeric-looking text (IDs /
# ranks / scores) that must stay stored as text (matches the original
# inlineStr cells) instead of being auto-coerced to numbers by the COM
# layer. Mark them as Text format before writing the values.
$ws.Range("A103:B111").NumberFormat = "@"
$ws.Range("E103:E111").NumberFormat = "@"
# C104's replacement value is also purely numeric-looking text.
$ws.Range("C104").NumberFormat = "@"

# Row 103
$ws.Range("A103").Value = "0"
$ws.Range("B103").Value = "6010122"
$ws.Range("C103").Value = '"Edward Peng"'
$ws.Range("E103").Value = "0"

# Row 104
$ws.Range("A104").Value = "0"
$ws.Range("B104").Value = "8850180"
$ws.Range("C104").Value = "30624300"
$ws.Range("E104").Value = "0"

# Row 105
$ws.Range("A105").Value = "0"
$ws.Range("B105").Value = "9195340"
$ws.Range("C105").Value = "Namllllllik"
$ws.Range("E105").Value = "0"

# Row 106
$ws.Range("A106").Value = "0"
$ws.Range("B106").Value = "9913517"
$ws.Range("C106").Value = '"Kenny Chan"'
$ws.Range("E106").Value = "0"

# Row 107
$ws.Range("A107").Value = "0"
$ws.Range("B107").Value = "10636651"
$ws.Range("C107").Value = '"Ismail Aflou"'
$ws.Range("E107").Value = "0"

# Row 108
$ws.Range("A108").Value = "0"
$ws.Range("B108").Value = "12648101"
$ws.Range("C108").Value = '"player 198827"'
$ws.Range("E108").Value = "0"

# Row 109
$ws.Range("A109").Value = "0"
$ws.Range("B109").Value = "15755724"
$ws.Range("C109").Value = '"Last Good"'
$ws.Range("E109").Value = "0"

# Row 110
$ws.Range("A110").Value = "0"
$ws.Range("B110").Value = "28624723"
$ws.Range("C110").Value = '"Woody Shade"'
$ws.Range("E110").Value = "0"

# Row 111
$ws.Range("A111").Value = "52311"
$ws.Range("B111").Value = "41848598"
$ws.Range("C111").Value = "国家一级保护沙雕"
$ws.Range("E111").Value = "3224"

# Remove the old rows 112-151 entirely so the sheet now ends at row 111.
$ws.Rows("112:151").Delete()
